$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data range A2:D18 by column A ascending (time column),
# re-establishing chronological order after needle calibration.
$rng = $ws.Range("A2:D18")
$sortKey = $ws.Range("A2:A18")

$rng.Sort($sortKey, 1)
